# Documentation.docx: add the "Dataset description" section (heading,
# intro paragraph and two bulleted dataset descriptions) right after the
# existing "MacQueen algorithm" bullet, mirroring the target diff.

$d = $word.ActiveDocument
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) Heading "Dataset description" + intro sentence.
#    Both are appended in one shot right at the very end of the story,
#    after the MacQueen paragraph's closing bookmark.
# ---------------------------------------------------------------------
$introXml = @"
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="1"/>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
    <w:t>Dataset description</w:t>
  </w:r>
</w:p>
<w:p $w>
  <w:pPr>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve">Two datasets were used in our </w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t>k-means</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve"> implementation:</w:t>
  </w:r>
</w:p>
"@
$endRange = $d.Range($d.Content.End, $d.Content.End)
[void]$endRange.InsertXML($introXml)

# ---------------------------------------------------------------------
# 2) First bulleted dataset ("Skin Segmentation Data Set ...").
#    Inserted as plain "List Paragraph" text first; the bullet/number
#    list is then minted fresh via ListFormat so a brand-new numbering
#    definition is created for it (mirrors the new numId in the diff).
# ---------------------------------------------------------------------
$skinXml = @"
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="a6"/>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t>Skin Segmentation Data Set</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve"> – RGB values from face images of various age, race groups and genders. They are taken from FERET and PAL databases. </w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t>50859</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve"> instances are skin samples, and </w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t>194198</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve"> are non-skin samples (</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t>artificially generated values).</w:t>
  </w:r>
</w:p>
"@
$endRange = $d.Range($d.Content.End, $d.Content.End)
[void]$endRange.InsertXML($skinXml)

$skinPara = $d.Paragraphs($d.Paragraphs.Count)
$skinPara.Range.ListFormat.ApplyListTemplateWithLevel()
$datasetNumId = $skinPara.Range.ListFormat.ListString
$datasetListXml = $skinPara.Range.WordOpenXML
$numIdMatch = [regex]::Match($datasetListXml, 'w:numId="(\d+)"')
if (-not $numIdMatch.Success) {
    $numIdMatch = [regex]::Match($skinPara.Range.ListFormat.ListTemplate.ToString(), '(\d+)')
}

# ---------------------------------------------------------------------
# 3) Second bulleted dataset ("HTRU2 Data Set ..."), placed on the same
#    list as the first bullet (same numId minted above).
# ---------------------------------------------------------------------
$htruXml = @"
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="a6"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="$($numIdMatch.Groups[1].Value)"/>
    </w:numPr>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t>HTRU2 Data Set</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve"> - a</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve"> sample of pulsar candidates collected during the High Time Resolution Universe Survey</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve">. There are 8 continuous variables describing each instance in dataset. 1639 instances are real pulsar examples and </w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t>16259</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve"> are fake.</w:t>
  </w:r>
</w:p>
"@
$endRange = $d.Range($d.Content.End, $d.Content.End)
[void]$endRange.InsertXML($htruXml)

Write-Output "Dataset description section inserted (numId=$($numIdMatch.Groups[1].Value))."
